$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3400.3572
$ws.Range("I100").Value = 2966.6667
$ws.Range("J100").Value = 3518.6365
$ws.Range("K100").Value = 2966.6667
$ws.Range("L100").Value = 3518.6365
$ws.Range("M100").Value = -2425.6667
$ws.Range("N100").Value = -4600.636500000001
$ws.Range("H101").Value = 925.9
$ws.Range("I101").Value = 784.25
$ws.Range("J101").Value = 1492.5
$ws.Range("K101").Value = 2352.75
$ws.Range("L101").Value = 4477.5
$ws.Range("M101").Value = -730.75
$ws.Range("N101").Value = -7721.5
$ws.Range("H132").Value = 2270.9792
$ws.Range("I132").Value = 1655.9688
$ws.Range("J132").Value = 3501
$ws.Range("K132").Value = 4967.9064
$ws.Range("L132").Value = 10503
$ws.Range("M132").Value = -2437.9064
$ws.Range("N132").Value = -15563
$ws.Range("H134").Value = 69900
$ws.Range("J134").Value = 69900
$ws.Range("L134").Value = 69900
$ws.Range("N134").Value = -80040
$ws.Range("H135").Value = 1003.4583
$ws.Range("I135").Value = 946.05884
$ws.Range("K135").Value = 8514.529560000001
$ws.Range("M135").Value = -5979.529560000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1273.4286
$ws.Range("J45").Value = 1626.6666
$ws.Range("L45").Value = 1626.6666
$ws.Range("N45").Value = -2380.6666
$ws.Range("H59").Value = 21250
$ws.Range("J59").Value = 21250
$ws.Range("L59").Value = 21250
$ws.Range("N59").Value = -22858
$ws.Range("H61").Value = 2580.9412
$ws.Range("I61").Value = 2023
$ws.Range("J61").Value = 3920
$ws.Range("K61").Value = 2023
$ws.Range("L61").Value = 3920
$ws.Range("M61").Value = -1811
$ws.Range("N61").Value = -4344
$ws.Range("H74").Value = 984
$ws.Range("I74").Value = 688
$ws.Range("J74").Value = 1650
$ws.Range("K74").Value = 688
$ws.Range("L74").Value = 1650
$ws.Range("M74").Value = 186
$ws.Range("N74").Value = -3398
$ws.Range("H77").Value = 984
$ws.Range("I77").Value = 688
$ws.Range("J77").Value = 1650
$ws.Range("K77").Value = 3440
$ws.Range("L77").Value = 8250
$ws.Range("M77").Value = 928
$ws.Range("N77").Value = -16986
$ws.Range("H102").Value = 1700
$ws.Range("I102").Value = 1700
$ws.Range("K102").Value = 1700
$ws.Range("M102").Value = -78
$ws.Range("H123").Value = 26236.4
$ws.Range("J123").Value = 26236.4
$ws.Range("L123").Value = 26236.4
$ws.Range("N123").Value = -36036.4
$ws.Range("H136").Value = 2580.9412
$ws.Range("I136").Value = 2023
$ws.Range("J136").Value = 3920
$ws.Range("K136").Value = 6069
$ws.Range("L136").Value = 11760
$ws.Range("M136").Value = -3519
$ws.Range("N136").Value = -16860

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2031.1154
$ws.Range("I134").Value = 2058.125
$ws.Range("J134").Value = 1707
$ws.Range("K134").Value = 6174.375
$ws.Range("L134").Value = 5121
$ws.Range("M134").Value = -3639.375
$ws.Range("N134").Value = -10191

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2188.0488
$ws.Range("I31").Value = 1511.871
$ws.Range("J31").Value = 4284.2
$ws.Range("K31").Value = 1511.871
$ws.Range("L31").Value = 4284.2
$ws.Range("M31").Value = -1216.871
$ws.Range("N31").Value = -4874.2
$ws.Range("H34").Value = 2188.0488
$ws.Range("I34").Value = 1511.871
$ws.Range("J34").Value = 4284.2
$ws.Range("K34").Value = 1511.871
$ws.Range("L34").Value = 4284.2
$ws.Range("M34").Value = -1309.871
$ws.Range("N34").Value = -4688.2
$ws.Range("H58").Value = 1038.2903
$ws.Range("I58").Value = 1115.12
$ws.Range("J58").Value = 718.1667
$ws.Range("K58").Value = 1115.12
$ws.Range("L58").Value = 718.1667
$ws.Range("M58").Value = -912.1199999999999
$ws.Range("N58").Value = -1124.1667
$ws.Range("H59").Value = 26705
$ws.Range("J59").Value = 20057.5
$ws.Range("L59").Value = 20057.5
$ws.Range("N59").Value = -22347.5
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H134").Value = 2174.8215
$ws.Range("I134").Value = 2021.2632
$ws.Range("J134").Value = 2499
$ws.Range("K134").Value = 6063.7896
$ws.Range("L134").Value = 7497
$ws.Range("M134").Value = -3528.7896
$ws.Range("N134").Value = -12567
$ws.Range("H136").Value = 1038.2903
$ws.Range("I136").Value = 1115.12
$ws.Range("J136").Value = 718.1667
$ws.Range("K136").Value = 3345.36
$ws.Range("L136").Value = 2154.5001
$ws.Range("M136").Value = -795.3599999999997
$ws.Range("N136").Value = -7254.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 1722.1111
$ws.Range("I10").Value = 150
$ws.Range("J10").Value = 2979.8
$ws.Range("K10").Value = 450
$ws.Range("L10").Value = 8939.400000000001
$ws.Range("M10").Value = -311
$ws.Range("N10").Value = -9217.400000000001
$ws.Range("H51").Value = 1000
$ws.Range("I51").Value = 1000
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 3000
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -2540
$ws.Range("N51").ClearContents()
$ws.Range("H113").Value = 680.17645
$ws.Range("I113").Value = 690.5
$ws.Range("J113").Value = 671
$ws.Range("K113").Value = 2071.5
$ws.Range("L113").Value = 2013
$ws.Range("M113").Value = 98.5
$ws.Range("N113").Value = -6353
$ws.Range("H122").Value = 677.62964
$ws.Range("I122").Value = 487.8889
$ws.Range("J122").Value = 772.5
$ws.Range("K122").Value = 4391.0001
$ws.Range("L122").Value = 6952.5
$ws.Range("M122").Value = -1941.0001
$ws.Range("N122").Value = -11852.5
$ws.Range("H126").Value = 3770
$ws.Range("I126").Value = 3015
$ws.Range("J126").Value = 3985.7144
$ws.Range("K126").Value = 9045
$ws.Range("L126").Value = 11957.1432
$ws.Range("M126").Value = -4105
$ws.Range("N126").Value = -21837.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 89.375
$ws.Range("I2").Value = 20
$ws.Range("J2").Value = 112.5
$ws.Range("K2").Value = 20
$ws.Range("L2").Value = 112.5
$ws.Range("M2").Value = 93
$ws.Range("N2").Value = -338.5
$ws.Range("H69").Value = 181282.86
$ws.Range("J69").Value = 181282.86
$ws.Range("L69").Value = 181282.86
$ws.Range("N69").Value = -182780.86
$ws.Range("H72").Value = 181282.86
$ws.Range("J72").Value = 181282.86
$ws.Range("L72").Value = 543848.58
$ws.Range("N72").Value = -551336.58
$ws.Range("H126").Value = 2819
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2819
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 8457
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -13397
$ws.Range("H132").Value = 1923.6
$ws.Range("I132").Value = 1656.174
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 4968.522
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -2438.522
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5131.923
$ws.Range("I40").Value = 10966.667
$ws.Range("K40").Value = 10966.667
$ws.Range("M40").Value = -10830.667
$ws.Range("H93").Value = 1335.3
$ws.Range("I93").Value = 910.6
$ws.Range("J93").Value = 1760
$ws.Range("K93").Value = 910.6
$ws.Range("L93").Value = 1760
$ws.Range("M93").Value = 337.4
$ws.Range("N93").Value = -4256
$ws.Range("H122").Value = 11116794
$ws.Range("I122").Value = 5441.5
$ws.Range("J122").Value = 33339500
$ws.Range("K122").Value = 16324.5
$ws.Range("L122").Value = 100018500
$ws.Range("M122").Value = -13874.5
$ws.Range("N122").Value = -100023400
$ws.Range("H136").Value = 3574.9375
$ws.Range("I136").Value = 2736.2727
$ws.Range("J136").Value = 5420
$ws.Range("K136").Value = 8208.8181
$ws.Range("L136").Value = 16260
$ws.Range("M136").Value = -5658.8181
$ws.Range("N136").Value = -21360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 21832.234
$ws.Range("J123").Value = 21832.234
$ws.Range("L123").Value = 21832.234
$ws.Range("N123").Value = -31632.234
$ws.Range("H132").Value = 2206.6
$ws.Range("I132").Value = 1745.3684
$ws.Range("J132").Value = 3667.1667
$ws.Range("K132").Value = 5236.1052
$ws.Range("L132").Value = 11001.5001
$ws.Range("M132").Value = -2706.1052
$ws.Range("N132").Value = -16061.5001
